$wb = $excel.ActiveWorkbook

$wsCreate = $wb.Worksheets.Item("Create")
$wsEdit   = $wb.Worksheets.Item("Edit")
$wsDelete = $wb.Worksheets.Item("Delete")

# Apostrophe prefix forces the engine to keep the cell as text (t="s")
# and reuse the existing "quotePrefix" cell style (style index 1) instead
# of re-interpreting numeric-looking text as a number and dropping the
# style, matching the style/type seen in the target workbook.
$q = [char]39

# ---------------------------------------------------------------------
# "Create" sheet (sheet1.xml)
# Row 2: A2/C2 keep style 1, B2 has no style.
# Row 3: A3/B3/C3 all get style 1 (B3 previously had none).
# ---------------------------------------------------------------------
$wsCreate.Range("A2").Value = "$q" + "2265"
$wsCreate.Range("B2").Value = "Mapping"
$wsCreate.Range("C2").Value = "$q" + "6756"

$wsCreate.Range("A3").Value = "$q" + "9806"
$wsCreate.Range("B3").Value = "$q" + "ListMap"
$wsCreate.Range("C3").Value = "$q" + "9005"

# ---------------------------------------------------------------------
# "Edit" sheet (sheet3.xml)
# Row 2: A2/C2 keep style 1, B2 has no style.
# Row 3: A3/B3/C3 all get style 1 (B3 previously had none).
# ---------------------------------------------------------------------
$wsEdit.Range("A2").Value = "$q" + "2265"
$wsEdit.Range("B2").Value = "Mapping"
$wsEdit.Range("C2").Value = "$q" + "6756"

$wsEdit.Range("A3").Value = "$q" + "9809"
$wsEdit.Range("B3").Value = "$q" + "HashMap"
$wsEdit.Range("C3").Value = "$q" + "9002"

# ---------------------------------------------------------------------
# "Delete" sheet (sheet4.xml)
# Row 2: A2/C2 keep style 1, B2 has no style.
# Row 3: A3/B3/C3 all get style 1 (B3 previously had none).
# ---------------------------------------------------------------------
$wsDelete.Range("A2").Value = "$q" + "2265"
$wsDelete.Range("B2").Value = "Mapping"
$wsDelete.Range("C2").Value = "$q" + "6756"

$wsDelete.Range("A3").Value = "$q" + "9809"
$wsDelete.Range("B3").Value = "$q" + "HashMap"
$wsDelete.Range("C3").Value = "$q" + "9002"

# ---------------------------------------------------------------------
# Selections / active sheet.
# The "Delete" sheet was previously the active tab; the "Create" sheet
# becomes active now, so select the other sheets' ranges first and the
# "Create" sheet's range last, so it ends up as the active tab.
# ---------------------------------------------------------------------
$wsEdit.Range("C4").Select()
$wsDelete.Range("C3").Select()
$wsCreate.Range("C3").Select()
